$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.903.58"
$ws.Range("E2").Value = "  +0.11%  "
$ws.Range("D3").Value = "1.589.77"
$ws.Range("E3").Value = "  -1.65%  "
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("E5").Value = "  -1.25%  "
$ws.Range("E6").Value = "  -0.17%  "
$ws.Range("E7").Value = "  -3.22%  "
$ws.Range("E8").Value = "  -0.43%  "
$ws.Range("E9").Value = "  +0.27%  "
$ws.Range("D10").Value = "18.26"
$ws.Range("E10").Value = "  -0.63%  "
$ws.Range("D11").Value = "0.0789"
$ws.Range("E11").Value = "  -0.30%  "
$ws.Range("D12").Value = "1.810.56"
$ws.Range("E12").Value = "  -1.63%  "
$ws.Range("D13").Value = "1.589.89"
$ws.Range("E13").Value = "  -1.61%  "
$ws.Range("E14").Value = "  -2.56%  "
$ws.Range("D15").Value = "0.513"
$ws.Range("E15").Value = "  -2.06%  "
$ws.Range("D16").Value = "25.911.62"
$ws.Range("D17").Value = "60.26"
$ws.Range("E17").Value = "  -1.92%  "
$ws.Range("E18").Value = "  -1.54%  "
$ws.Range("E19").Value = "  -0.22%  "
$ws.Range("D20").Value = "193.76"
$ws.Range("E20").Value = "  +1.40%  "
$ws.Range("E21").Value = "  -0.73%  "
$ws.Range("D22").Value = "9.41"
$ws.Range("E22").Value = "  -0.85%  "
$ws.Range("D23").Value = "5.94"
$ws.Range("E23").Value = "  -1.21%  "
$ws.Range("E24").Value = "  -1.61%  "
$ws.Range("E25").Value = "  -1.48%  "
$ws.Range("E26").Value = "  -0.22%  "
$ws.Range("E27").Value = "  -0.19%  "
$ws.Range("D28").Value = "15.15"
$ws.Range("E28").Value = "  -0.49%  "
$ws.Range("E29").Value = "  -2.39%  "
$ws.Range("E30").Value = "  -5.36%  "
$ws.Range("D31").Value = "0.0473"
$ws.Range("E31").Value = "  -0.68%  "
$ws.Range("E32").Value = "  +0.01%  "
$ws.Range("E33").Value = "  -1.62%  "
$ws.Range("E34").Value = "  +1.12%  "
$ws.Range("E35").Value = "  -2.23%  "
$ws.Range("D36").Value = "1.108.01"
$ws.Range("E36").Value = "  -1.53%  "
$ws.Range("E37").Value = "  -0.28%  "
$ws.Range("E38").Value = "  -1.43%  "
$ws.Range("E39").Value = "  -0.77%  "
$ws.Range("E40").Value = "  -1.54%  "
$ws.Range("E41").Value = "  -6.18%  "
$ws.Range("D42").Value = "0.819"
$ws.Range("E42").Value = "  +9.09%  "
$ws.Range("E43").Value = "  +2.40%  "
$ws.Range("D44").Value = "93.63"
$ws.Range("E44").Value = "  -4.63%  "
$ws.Range("D45").Value = "1.723.36"
$ws.Range("E45").Value = "  -1.64%  "
$ws.Range("D46").Value = "0.0₆0112"
$ws.Range("E46").Value = "  -0.59%  "
$ws.Range("D47").Value = "1.51"
$ws.Range("E47").Value = "  +0.03%  "
$ws.Range("D48").Value = "53.52"
$ws.Range("E48").Value = "  -0.77%  "
$ws.Range("E49").Value = "  -1.65%  "
$ws.Range("D50").Value = "0.408"
$ws.Range("E50").Value = "  -0.61%  "
$ws.Range("E51").Value = "  -0.20%  "
